$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# STREAMLINE sheet (sheet12 / "STREAMLINE"): add the 3rd image slot
# (SECTION_2_IMG_3) to the "Why choosing Streamline?" section, and the
# matching IMAGE::STREAMLINE_x references for the 3 images.
# ---------------------------------------------------------------------------
$wsStreamline = $wb.Worksheets.Item("STREAMLINE")

# Insert the 3 new data rows (top to bottom, each insert shifts the rows
# below it down by one, so the target row index is the same for all three
# since they land in order 12, 15, 18 of the *evolving* sheet).
$wsStreamline.Rows.Item(12).Insert()
$wsStreamline.Rows.Item(15).Insert()
$wsStreamline.Rows.Item(18).Insert()

# Row 12: SECTION_2_IMG_1 -> IMAGE::STREAMLINE_1
$wsStreamline.Cells.Item(12, 1).Value = "SECTION_2_IMG_1"

# Row 15: SECTION_2_IMG_2 -> IMAGE::STREAMLINE_2
$wsStreamline.Cells.Item(15, 1).Value = "SECTION_2_IMG_2"

# Row 18: SECTION_2_IMG_3 (brand new id) -> IMAGE::STREAMLINE_3
$wsStreamline.Cells.Item(18, 1).Value = "SECTION_2_IMG_3"

# ---------------------------------------------------------------------------
# IMAGE sheet (sheet3): register the 3 new streamline images
# ---------------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("IMAGE")

$wsImage.Cells.Item(95, 1).Value = "STREAMLINE_1"
$wsImage.Cells.Item(96, 1).Value = "STREAMLINE_2"
$wsImage.Cells.Item(97, 1).Value = "STREAMLINE_3"

$wsImage.Cells.Item(95, 2).Value = "/images/streamline/time.png"
$wsImage.Cells.Item(96, 2).Value = "/images/streamline/task.png"
$wsImage.Cells.Item(97, 2).Value = "/images/streamline/money.png"

# ---------------------------------------------------------------------------
# Back on STREAMLINE: fill in the IMAGE:: references, in the same order the
# original author entered them (2, 3, 1) so shared-string indices line up.
# ---------------------------------------------------------------------------
$wsStreamline.Cells.Item(15, 2).Value = "IMAGE::STREAMLINE_2"
$wsStreamline.Cells.Item(18, 2).Value = "IMAGE::STREAMLINE_3"
$wsStreamline.Cells.Item(12, 2).Value = "IMAGE::STREAMLINE_1"

# Row-height tweaks that came along with the new rows.
$wsStreamline.Rows.Item(14).RowHeight = 45
$wsStreamline.Rows.Item(17).RowHeight = 60
$wsStreamline.Rows.Item(29).RowHeight = 30

# ---------------------------------------------------------------------------
# Restore the view/selection state recorded in the saved file.
# ---------------------------------------------------------------------------
$wsImage.Activate()
$wsImage.Range("A96").Select()

$wsStreamline.Activate()
$wsStreamline.Range("B13").Select()
